$wb = $excel.ActiveWorkbook

# Sheet "SQL" - row 15
$ws = $wb.Worksheets.Item("SQL")
$ws.Range("B15").Value = "TestTestTestTestTestTestTestadasd"
$ws.Range("C15").Value = "TestTestTestTestTestTestTestsadasd"
$ws.Range("D15").Value = "TestTestTestTestTestTestTestasdsad"

# Sheet "Python" - row 32
$ws = $wb.Worksheets.Item("Python")
$ws.Range("B32").Value = "пвавапавпавasdasd"
$ws.Range("C32").Value = "павпвапавпвапasdasd"
$ws.Range("D32").Value = "вапвввввввввasdad"

# Sheet "Links" - row 1
$ws = $wb.Worksheets.Item("Links")
$ws.Range("B1").Value = "TestTestTestTestTestTestTestgdgdfasd"
$ws.Range("C1").Value = "TestTestTestTestTestTestTestdfgdfgasdasd"
$ws.Range("D1").Value = "TestTestTestTestTestTestTestdfgdfasdas"

# Sheet "Bash" - row 67
$ws = $wb.Worksheets.Item("Bash")
$ws.Range("B67").Value = "fsdfffffffffffffffffdasdsafsdfffffffffffffffffasdasdfsdfffffffffffffffffasdasd"
$ws.Range("C67").Value = "fsdfffffffffffffffffasdasdfsdfffffffffffffffffasdasdfsdfffffffffffffffffasdasd"
